# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Re-sequences the worker/period rows 16-25 on "Hoja1": Edinson Alberto
# Camacho Camacho's two periods move to the top (1903 before 1904), Leydis
# Diaz Palomo's row moves down, and Heidys Gonzalez Causil's four periods
# are interleaved with Mileidis Hueto Montalban's single period (ordered
# 2301, 2301, 2302, 2303, 2304) before Javier F Guerra Guerra and Lisbeth
# Paola Perez Padilla close out the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: EDINSON ALBERTO CAMACHO CAMACHO - periodo 1903
$ws.Range("C16").Value = "7917927"
$ws.Range("D16").Value = "EDINSON ALBERTO CAMACHO CAMACHO"
$ws.Range("E16").Value = "1903"
$ws.Range("F16").Value = 80000
$ws.Range("G16").Value = 2000000

# Row 17: EDINSON ALBERTO CAMACHO CAMACHO - periodo 1904
$ws.Range("C17").Value = "7917927"
$ws.Range("D17").Value = "EDINSON ALBERTO CAMACHO CAMACHO"
$ws.Range("E17").Value = "1904"
$ws.Range("F17").Value = 80000
$ws.Range("G17").Value = 2000000

# Row 18: LEYDIS DIAZ PALOMO - periodo 2211
$ws.Range("C18").Value = "1047367581"
$ws.Range("D18").Value = "LEYDIS DIAZ PALOMO"
$ws.Range("E18").Value = "2211"
$ws.Range("F18").Value = 1600
$ws.Range("G18").Value = 1200000

# Row 19: HEIDYS GONZALEZ CAUSIL - periodo 2301
$ws.Range("C19").Value = "1047467797"
$ws.Range("D19").Value = "HEIDYS GONZALEZ CAUSIL"
$ws.Range("E19").Value = "2301"
$ws.Range("F19").Value = 23200
$ws.Range("G19").Value = 1160000

# Row 20: MILEIDIS HUETO MONTALBAN - periodo 2301
$ws.Range("C20").Value = "1050962245"
$ws.Range("D20").Value = "MILEIDIS HUETO MONTALBAN"
$ws.Range("E20").Value = "2301"
$ws.Range("F20").Value = 34800
$ws.Range("G20").Value = 870000

# Row 21: HEIDYS GONZALEZ CAUSIL - periodo 2302
$ws.Range("C21").Value = "1047467797"
$ws.Range("D21").Value = "HEIDYS GONZALEZ CAUSIL"
$ws.Range("E21").Value = "2302"
$ws.Range("F21").Value = 46400
$ws.Range("G21").Value = 1160000

# Row 22: HEIDYS GONZALEZ CAUSIL - periodo 2303
$ws.Range("C22").Value = "1047467797"
$ws.Range("D22").Value = "HEIDYS GONZALEZ CAUSIL"
$ws.Range("E22").Value = "2303"
$ws.Range("F22").Value = 46400
$ws.Range("G22").Value = 1160000

# Row 23: HEIDYS GONZALEZ CAUSIL - periodo 2304
$ws.Range("C23").Value = "1047467797"
$ws.Range("D23").Value = "HEIDYS GONZALEZ CAUSIL"
$ws.Range("E23").Value = "2304"
$ws.Range("F23").Value = 37120
$ws.Range("G23").Value = 1160000

# Row 24: JAVIER F GUERRA GUERRA - periodo 2311
$ws.Range("C24").Value = "73108064"
$ws.Range("D24").Value = "JAVIER F GUERRA GUERRA"
$ws.Range("E24").Value = "2311"
$ws.Range("F24").Value = 46400
$ws.Range("G24").Value = 1160000

# Row 25: LISBETH PAOLA PEREZ PADILLA - periodo 2311 (unchanged values, kept for completeness)
$ws.Range("C25").Value = "1002199218"
$ws.Range("D25").Value = "LISBETH PAOLA PEREZ PADILLA"
$ws.Range("E25").Value = "2311"
$ws.Range("F25").Value = 46400
$ws.Range("G25").Value = 1160000

Write-Output "rows 16-25 updated"
